# Adds the "ExposureLoops" feature described in the commit message:
#  - a new italic note under the first ("For [y] = 4 and [z] = i") table
#    explaining the new three-point UV alignment settings
#  - a new row (#14) in that same table describing the new ExposureLoops
#    variable (int) used to specify the number of exposure loops

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S-Matrix")

# New italic annotation next to the "For [y] = 4 and [z] = i" header (row 4)
$ws.Range("C4").Value = "Three-point UV alignment settings"
$ws.Range("C4").Font.Italic = $true

# New table row describing the ExposureLoops variable
$ws.Range("A31").Value = 14
$ws.Range("B31").Value = "ExposureLoops"
$ws.Range("C31").Value = "Number of exposure loops"
$ws.Range("D31").Value = "int"
